$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.459666967391968
$ws.Range("B1").Value = 3.819129467010498
$ws.Range("C1").Value = 2.783946514129639
$ws.Range("D1").Value = 0.891598105430603
$ws.Range("E1").Value = 1.091726660728455
